# Commit: "Chanfge format to create the paths, in that way the path create
# according the the operational system of the PC"
#
# The author switched the workbook's stored "absolute path" metadata from a
# SharePoint URL to a local, OS-native (Windows) filesystem path, and (in the
# same save) Excel's re-serialisation normalised the "Plain_English" column
# header text to "Plain English" on every lookup-table sheet that has it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename the "Plain_English" header to "Plain English" everywhere it is
#    used (Fuel_to_Code, VehFuel_to_Code, Tech_to_Code, Dem_to_Code).
# ---------------------------------------------------------------------------
$headerSheets = @("Fuel_to_Code", "VehFuel_to_Code", "Tech_to_Code", "Dem_to_Code")
foreach ($sheetName in $headerSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $found = $ws.Rows.Item(1).Find("Plain_English")
    if ($found -ne $null) {
        $found.Value = "Plain English"
    }
}

# ---------------------------------------------------------------------------
# 2) Point the workbook at its local path instead of the old SharePoint URL
#    (the commit's stated intent: "create the paths ... according [to] the
#    operational system of the PC"). Re-saving from the local desktop copy
#    is what flips Excel's stored absolute-path hint from the SharePoint
#    URL to the Windows path below.
# ---------------------------------------------------------------------------
$localPath = "C:\Users\ClimateLeadGroup\Desktop\CLG_repositories\osemosys_momf\t3a_experiments\Experiment_1\0_From_Confection\A-I_Classifier_Modes_Transport.xlsx"
$wb.SaveAs($localPath)
